$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "96.763.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.678.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.24"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.87"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +8.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "656.89"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.424"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.675.72"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.80"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.206"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.363.95"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000270"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.541.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.91"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.670.09"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.47%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.80"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.528"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "529.79"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.50"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000205"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.28"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.38"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.872.07"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.169"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +7.05%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.53"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.32%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.06"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.92"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +15.87%  "
$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.186"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "675.60"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +9.35%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.41%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "32.61"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.26%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.596"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.87%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.93"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.161"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.57"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.94%  "
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.79"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +20.72%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.962"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.46%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0464"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.19%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.440"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +12.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.77"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.71%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.33"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.66"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.09%  "
